$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.207.10'
$ws.Range('E2').Value = '  +2.88%  '
$ws.Range('D3').Value = '2.306.22'
$ws.Range('E3').Value = '  +1.87%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.02'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.19'
$ws.Range('E6').Value = '  +5.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.537'
$ws.Range('E7').Value = '  +2.02%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.525'
$ws.Range('E9').Value = '  +6.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.00'
$ws.Range('E10').Value = '  +2.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0818'
$ws.Range('E11').Value = '  +3.30%  '
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('E13').Value = '  +5.52%  '
$ws.Range('D14').Value = '2.663.84'
$ws.Range('E14').Value = '  +1.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.98'
$ws.Range('E15').Value = '  +4.04%  '
$ws.Range('D16').Value = '2.313.09'
$ws.Range('E16').Value = '  +2.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.812'
$ws.Range('E17').Value = '  +2.63%  '
$ws.Range('D18').Value = '43.118.11'
$ws.Range('E18').Value = '  +3.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.53'
$ws.Range('E19').Value = '  +1.00%  '
$ws.Range('D20').Value = '0.0₃0920'
$ws.Range('E20').Value = '  +2.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.10'
$ws.Range('E21').Value = '  +2.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.39'
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.92'
$ws.Range('E23').Value = '  +1.50%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.63'
$ws.Range('E24').Value = '  +2.89%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.01'
$ws.Range('E25').Value = '  +4.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.68'
$ws.Range('E27').Value = '  +4.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.46'
$ws.Range('E28').Value = '  +2.17%  '
$ws.Range('E29').Value = '  +2.10%  '
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '166.71'
$ws.Range('E31').Value = '  +4.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.32'
$ws.Range('E32').Value = '  +2.20%  '
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.12'
$ws.Range('E34').Value = '  -1.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.72'
$ws.Range('E35').Value = '  +3.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0740'
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range('E37').Value = '  +2.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.39'
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('E39').Value = '  +1.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.83'
$ws.Range('E40').Value = '  +0.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.30'
$ws.Range('E41').Value = '  +7.67%  '
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').Value = '1.986.75'
$ws.Range('E43').Value = '  +0.77%  '
$ws.Range('E44').Value = '  +2.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.12'
$ws.Range('E45').Value = '  +2.59%  '
$ws.Range('E46').Value = '  +3.15%  '
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('E48').Value = '  +18.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.76'
$ws.Range('E49').Value = '  +5.37%  '
$ws.Range('D50').Value = '2.532.54'
$ws.Range('E50').Value = '  +1.89%  '
$ws.Range('E51').Value = '  +1.95%  '
